$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Hydrogen): corrected Iron & steel hydrogen demand, clear the
#     now-unused Non-metallic minerals figure ---
$ws.Range("B3").Value = 9206820.142651889
$ws.Range("D3").ClearContents()

# --- Row 4 (Methanol): corrected Chemicals figure ---
$ws.Range("C4").Value = 9.768018748098955

# --- Row 5 (Ammonia): corrected Chemicals figure ---
$ws.Range("C5").Value = 1740.714904462766

# --- Row 7: "Other" is now specifically "Biogas", with a corrected value ---
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 1372.1853523845

# --- New row 8: re-introduce "Other" below Biogas, carrying the formatting
#     of the row above it (bold label, bordered, centered) ---
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 1307.949500168818
